# Fix the algorithm/conditions on filtering the status of candidates.
# The Candidate / Status / Action Date values for rows 3, 5, 6, 7 need to be
# rotated: row7 -> row3, row3 -> row5, row5 -> row6, row6 -> row7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values (columns D, E, F) before overwriting anything.
# Use Value2 since plain .Value is not reliably readable in this runtime.
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2

$d5 = $ws.Range("D5").Value2
$e5 = $ws.Range("E5").Value2
$f5 = $ws.Range("F5").Value2

$d6 = $ws.Range("D6").Value2
$e6 = $ws.Range("E6").Value2
$f6 = $ws.Range("F6").Value2

$d7 = $ws.Range("D7").Value2
$e7 = $ws.Range("E7").Value2
$f7 = $ws.Range("F7").Value2

# Row 3 gets old row 7's values
$ws.Range("D3").Value2 = $d7
$ws.Range("E3").Value2 = $e7
$ws.Range("F3").Value2 = $f7

# Row 5 gets old row 3's values
$ws.Range("D5").Value2 = $d3
$ws.Range("E5").Value2 = $e3
$ws.Range("F5").Value2 = $f3

# Row 6 gets old row 5's values
$ws.Range("D6").Value2 = $d5
$ws.Range("E6").Value2 = $e5
$ws.Range("F6").Value2 = $f5

# Row 7 gets old row 6's values
$ws.Range("D7").Value2 = $d6
$ws.Range("E7").Value2 = $e6
$ws.Range("F7").Value2 = $f6
